# Groopman et al 2019 - Genes with diagnostic variants
# Commit: "Refined metadata to be additional tab"
#
# This script:
#  1. Updates the "time_taken" timestamps (column F, rows 2-67) on the
#     "data" sheet to reflect the refined/re-run query time.
#  2. Adds a new "metadata" worksheet (placed after "data") carrying the
#     panel-level metadata (name, id, version, version_created,
#     query_time, get_request) that used to live only implicitly.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

$dataSheet.Range("F2").Value = "2021-10-05 14:20:34.503774"
$dataSheet.Range("F3").Value = "2021-10-05 14:20:34.503782"
$dataSheet.Range("F4").Value = "2021-10-05 14:20:34.503785"
$dataSheet.Range("F5").Value = "2021-10-05 14:20:34.503788"
$dataSheet.Range("F6").Value = "2021-10-05 14:20:34.503791"
$dataSheet.Range("F7").Value = "2021-10-05 14:20:34.503793"
$dataSheet.Range("F8").Value = "2021-10-05 14:20:34.503796"
$dataSheet.Range("F9").Value = "2021-10-05 14:20:34.503798"
$dataSheet.Range("F10").Value = "2021-10-05 14:20:34.503801"
$dataSheet.Range("F11").Value = "2021-10-05 14:20:34.503804"
$dataSheet.Range("F12").Value = "2021-10-05 14:20:34.503807"
$dataSheet.Range("F13").Value = "2021-10-05 14:20:34.503809"
$dataSheet.Range("F14").Value = "2021-10-05 14:20:34.503812"
$dataSheet.Range("F15").Value = "2021-10-05 14:20:34.503814"
$dataSheet.Range("F16").Value = "2021-10-05 14:20:34.503817"
$dataSheet.Range("F17").Value = "2021-10-05 14:20:34.503819"
$dataSheet.Range("F18").Value = "2021-10-05 14:20:34.503822"
$dataSheet.Range("F19").Value = "2021-10-05 14:20:34.503825"
$dataSheet.Range("F20").Value = "2021-10-05 14:20:34.503827"
$dataSheet.Range("F21").Value = "2021-10-05 14:20:34.503830"
$dataSheet.Range("F22").Value = "2021-10-05 14:20:34.503832"
$dataSheet.Range("F23").Value = "2021-10-05 14:20:34.503835"
$dataSheet.Range("F24").Value = "2021-10-05 14:20:34.503838"
$dataSheet.Range("F25").Value = "2021-10-05 14:20:34.503840"
$dataSheet.Range("F26").Value = "2021-10-05 14:20:34.503843"
$dataSheet.Range("F27").Value = "2021-10-05 14:20:34.503846"
$dataSheet.Range("F28").Value = "2021-10-05 14:20:34.503848"
$dataSheet.Range("F29").Value = "2021-10-05 14:20:34.503851"
$dataSheet.Range("F30").Value = "2021-10-05 14:20:34.503853"
$dataSheet.Range("F31").Value = "2021-10-05 14:20:34.503856"
$dataSheet.Range("F32").Value = "2021-10-05 14:20:34.503858"
$dataSheet.Range("F33").Value = "2021-10-05 14:20:34.503861"
$dataSheet.Range("F34").Value = "2021-10-05 14:20:34.503864"
$dataSheet.Range("F35").Value = "2021-10-05 14:20:34.503866"
$dataSheet.Range("F36").Value = "2021-10-05 14:20:34.503869"
$dataSheet.Range("F37").Value = "2021-10-05 14:20:34.503871"
$dataSheet.Range("F38").Value = "2021-10-05 14:20:34.503874"
$dataSheet.Range("F39").Value = "2021-10-05 14:20:34.503876"
$dataSheet.Range("F40").Value = "2021-10-05 14:20:34.503879"
$dataSheet.Range("F41").Value = "2021-10-05 14:20:34.503882"
$dataSheet.Range("F42").Value = "2021-10-05 14:20:34.503885"
$dataSheet.Range("F43").Value = "2021-10-05 14:20:34.503887"
$dataSheet.Range("F44").Value = "2021-10-05 14:20:34.503890"
$dataSheet.Range("F45").Value = "2021-10-05 14:20:34.503892"
$dataSheet.Range("F46").Value = "2021-10-05 14:20:34.503895"
$dataSheet.Range("F47").Value = "2021-10-05 14:20:34.503897"
$dataSheet.Range("F48").Value = "2021-10-05 14:20:34.503900"
$dataSheet.Range("F49").Value = "2021-10-05 14:20:34.503902"
$dataSheet.Range("F50").Value = "2021-10-05 14:20:34.503905"
$dataSheet.Range("F51").Value = "2021-10-05 14:20:34.503907"
$dataSheet.Range("F52").Value = "2021-10-05 14:20:34.503910"
$dataSheet.Range("F53").Value = "2021-10-05 14:20:34.503912"
$dataSheet.Range("F54").Value = "2021-10-05 14:20:34.503915"
$dataSheet.Range("F55").Value = "2021-10-05 14:20:34.503918"
$dataSheet.Range("F56").Value = "2021-10-05 14:20:34.503920"
$dataSheet.Range("F57").Value = "2021-10-05 14:20:34.503923"
$dataSheet.Range("F58").Value = "2021-10-05 14:20:34.503926"
$dataSheet.Range("F59").Value = "2021-10-05 14:20:34.503928"
$dataSheet.Range("F60").Value = "2021-10-05 14:20:34.503931"
$dataSheet.Range("F61").Value = "2021-10-05 14:20:34.503933"
$dataSheet.Range("F62").Value = "2021-10-05 14:20:34.503936"
$dataSheet.Range("F63").Value = "2021-10-05 14:20:34.503938"
$dataSheet.Range("F64").Value = "2021-10-05 14:20:34.503941"
$dataSheet.Range("F65").Value = "2021-10-05 14:20:34.503943"
$dataSheet.Range("F66").Value = "2021-10-05 14:20:34.503947"
$dataSheet.Range("F67").Value = "2021-10-05 14:20:34.503950"

# --- Add the "metadata" worksheet, positioned right after "data" ---
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (B1:G1)
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Match the header formatting used on the "data" sheet (bold/bordered/
# centered style already present in the workbook) by copying it over,
# instead of constructing a brand-new style.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

# Data row (A2:G2)
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Groopman et al 2019 - Genes with diagnostic variants"
$metaSheet.Range("C2").Value = 720
$metaSheet.Range("D2").Value = "'0.8"
$metaSheet.Range("E2").Value = "2019-07-09T15:48:14.145108Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:20:34.500243"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/720/?format=json"

# Match A2's formatting (index-column style) to the "data" sheet as well.
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# Keep "data" as the active sheet/selection, matching the original
# workbook view state (only the sheet list itself changed).
$dataSheet.Activate() | Out-Null
$dataSheet.Range("A1").Select() | Out-Null
